# Improved snippets for https://github.com/2sic/2sxc/issues/976
#
# Adds 8 new "script/css optimization" snippet rows right after the
# existing "script" / "css, style-sheet" rows (242/243), turning the old
# single-purpose rows into a small family of related snippets, and shifts
# everything that used to follow them down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Make room: insert 6 blank rows right after row 243 (before old row 244).
#    This pushes old rows 244..260 down to 250..266 and keeps column-A's
#    quote-prefix style (s="1") because it is inherited from the row above.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows.Item(244).Insert()
}

# 2) Update the two existing rows (242 & 243) that stay in place but now
#    describe the "default" variant of the new mini-family, plus gain
#    F (description) and G (api-docs link) values.
$ws.Range("C242").Value = "script with optimization"
$ws.Range("F242").Value = "a script tag which enables scripts-combinations / packing etc. at default position & priority"
$ws.Range("G242").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("C243").Value = "script optimized into head"
$ws.Range("E243").Value = '<script src="@App.Path/dist/${1:myscripts}.js" type="text/javascript" data-enableoptimizations="500:head"></script>'
$ws.Range("F243").Value = "a script tag which enables scripts-combinations / packing etc. at lower priority in the header"
$ws.Range("G243").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

# 3) Fill in the 6 brand-new rows (244..249).
#    Column A gets a leading apostrophe so Excel stores it as literal text
#    with the quote-prefix cell style (s="1"), matching column A's style
#    on every other "@..." row in this table (e.g. rows 242/243 above it).
$ws.Range("A244").Value = "'@Html"
$ws.Range("B244").Value = "Resources"
$ws.Range("C244").Value = "script optimized into bottom of page"
$ws.Range("E244").Value = '<script src="@App.Path/dist/${1:myscripts}.js" type="text/javascript" data-enableoptimizations="bottom"></script>'
$ws.Range("F244").Value = "a script tag which enables scripts-combinations / packing etc. at default priority at the bottom"
$ws.Range("G244").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("A245").Value = "'@Html"
$ws.Range("B245").Value = "Resources"
$ws.Range("C245").Value = "script specifically in body"
$ws.Range("E245").Value = '<script src="@App.Path/dist/${1:myscripts}.js" type="text/javascript" data-enableoptimizations="175:body"></script>'
$ws.Range("F245").Value = "a script tag which enables scripts-combinations / packing etc. at custom priority at page-top (not header)"
$ws.Range("G245").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("A246").Value = "'@Html"
$ws.Range("B246").Value = "Resources"
$ws.Range("C246").Value = "css, style-sheet with optimization"
$ws.Range("E246").Value = '<link rel="stylesheet" href="@App.Path/dist/AppCatalog.css" data-enableoptimizations="100"/>'
$ws.Range("F246").Value = "a css-tag which enables optimizations - with priority 100 (default) in def. position"
$ws.Range("G246").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("A247").Value = "'@Html"
$ws.Range("B247").Value = "Resources"
$ws.Range("C247").Value = "css, style-sheet loaded in body"
$ws.Range("E247").Value = '<link rel="stylesheet" href="@App.Path/dist/AppCatalog.css" data-enableoptimizations="body"/>'
$ws.Range("F247").Value = "a css-tag which enables optimizations - with priority default inside body"
$ws.Range("G247").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("A248").Value = "'@Html"
$ws.Range("B248").Value = "Resources"
$ws.Range("C248").Value = "css, style-sheet loaded in head"
$ws.Range("E248").Value = '<link rel="stylesheet" href="@App.Path/dist/AppCatalog.css" data-enableoptimizations="150:head"/>'
$ws.Range("F248").Value = "a css-tag which enables optimizations - with custom priority in head"
$ws.Range("G248").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

$ws.Range("A249").Value = "'@Html"
$ws.Range("B249").Value = "Resources"
$ws.Range("C249").Value = "css, style-sheet loaded at bottom of page"
$ws.Range("E249").Value = '<link rel="stylesheet" href="@App.Path/dist/AppCatalog.css" data-enableoptimizations="bottom"/>'
$ws.Range("F249").Value = "a css-tag which enables optimizations - with def. priority at page bottom"
$ws.Range("G249").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

# 4) Grow the table (Table1) so it covers the 6 new rows too
#    (ref/autoFilter A1:G260 -> A1:G266).
$lo.Resize($ws.Range("A1:G266"))

# 5) Match the author's final viewport/selection in the sheet view.
$ws.Application.ActiveWindow.ScrollRow = 224
$ws.Range("C250").Select()
